$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scaling input changed (this cascades through the formulas for C4, E8:E11, E13:E15, E17)
$ws.Range("B4").Value = 40

# Row 8 - Maris Otter
$ws.Range("F8").Value = 11
$ws.Range("H8").Value = 110

# Row 9 - Crystal EBC 20
$ws.Range("F9").Value = 1.5
$ws.Range("H9").Value = 15

# Row 10 - Munich EBC 25
$ws.Range("F10").Value = 1.5
$ws.Range("H10").Value = 15

# Row 11 - Hvetemalt
$ws.Range("F11").Value = 1.2
$ws.Range("H11").Value = 12

# Row 13 - Chinook
$ws.Range("F13").Value = 100
$ws.Range("H13").Value = "(1pk)"

# Row 15 - Simcoe
$ws.Range("F15").Value = 100
$ws.Range("H15").Value = "(1 pk)"

# Row 14 - Cascade
$ws.Range("F14").Value = 400
$ws.Range("H14").Value = "(4 pk)"

# Row 17 - WLP 007 yeast
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = "(evt 2 hvis starter)"

# Update selection / active cell to match the saved workbook view
$ws.Range("R35").Select()
